$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (names re-ordered, total_registros updated),
# sorted descending by total_registros. Header row 1 and trailing row 12
# (71050834 / 1) are unchanged.
$data = @(
    @("PEREZ VEGA ANA YSABEL", 88),
    @("ZAPATA ZETA ROSA ARACELI", 81),
    @("GARAVITO LEON IVONNE LISSETH", 79),
    @("TIMOTEO BAYONA SHARYN LISSETH", 78),
    @("PANTA MONZON SHIRLEY MARIBEL", 75),
    @("NIÑO GUERRERO ANYELA MELINA", 70),
    @("CASTRO JUAREZ MARIA ISABEL", 66),
    @("VALLE SILVA SUTMMER ORFELINDA", 60),
    @("TIZON NUÑEZ FRESIA YAMILI", 57),
    @("CHERO JUAREZ ANYELA TATIANA", 38)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
